$wb = $excel.ActiveWorkbook

# =====================================================================
# 1. Insert a new "2022-Q1" fund-holding detail sheet, positioned right
#    before the "总计" (totals) summary sheet.
# =====================================================================
$totalSheet = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

# Reuse formatting from the most recent existing detail sheet ("2021-Q4")
# so the new sheet's header row / index column look the same.
$template = $wb.Worksheets.Item("2021-Q4")
$template.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$template.Range("A2").Copy()
$newSheet.Range("A2:A26").PasteSpecial(-4122)

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$data = @(
  @('005453', '前海开源医疗健康灵活配置混合A', '24.34', '89.12', '5.06', '1.2316', 9),
  @('005505', '前海开源中药研究精选股票A', '11.44', '91.95', '6.26', '0.7161', 10),
  @('005454', '前海开源医疗健康灵活配置混合C', '12.66', '89.12', '5.06', '0.6406', 9),
  @('000339', '长城医疗保健混合', '11.06', '88.40', '3.11', '0.3440', 10),
  @('005506', '前海开源中药研究精选股票C', '4.62', '91.95', '6.26', '0.2892', 10),
  @('011673', '长城医药科技六个月持有期混合型证券投资基金A', '8.65', '88.50', '3.12', '0.2699', 9),
  @('010054', '万家健康产业混合A', '8.13', '86.63', '3.17', '0.2577', 7),
  @('200006', '长城消费增值混合', '7.35', '88.38', '3.43', '0.2521', 7),
  @('001558', '天弘医疗健康混合A', '7.72', '76.93', '3.06', '0.2362', 9),
  @('519673', '银河康乐股票', '2.31', '92.35', '7.98', '0.1843', 1),
  @('001559', '天弘医疗健康混合C', '4.61', '76.93', '3.06', '0.1411', 9),
  @('008786', '长城健康生活灵活配置混合', '6.36', '78.57', '2.17', '0.1380', 5),
  @('011601', '前海开源公共卫生主题精选股票A', '2.25', '88.30', '5.07', '0.1141', 9),
  @('010055', '万家健康产业混合C', '3.36', '86.63', '3.17', '0.1065', 7),
  @('011602', '前海开源公共卫生主题精选股票C', '0.86', '88.30', '5.07', '0.0436', 9),
  @('010434', '红土创新医疗保健股票', '0.75', '92.96', '4.29', '0.0322', 8),
  @('200016', '长城稳健成长灵活配置混合', '0.83', '78.39', '3.53', '0.0293', 9),
  @('217021', '招商优势企业混合', '0.36', '69.72', '7.16', '0.0258', 1),
  @('011674', '长城医药科技六个月持有期混合型证券投资基金C', '0.72', '88.50', '3.12', '0.0225', 9),
  @('007254', '广发均衡价值混合', '0.49', '89.66', '4.43', '0.0217', 9),
  @('000649', '长城久鑫灵活配置混合', '0.46', '81.08', '3.20', '0.0147', 7),
  @('007381', '国融融信消费严选混合A', '0.27', '68.97', '2.70', '0.0073', 8),
  @('013072', '泰信医疗服务混合A', '0.10', '73.32', '3.16', '0.0032', 10),
  @('007382', '国融融信消费严选混合C', '0.04', '68.97', '2.70', '0.0011', 8),
  @('013073', '泰信医疗服务混合C', '0.00', '73.32', '3.16', '0', 10)
)

# Columns B (fund code) and D-G (fund scale/position/value) look numeric
# but must stay TEXT, matching the source workbook's convention. Force a
# text number format before assigning so Excel doesn't coerce them.
$newSheet.Range("B2:B26").NumberFormat = "@"
$newSheet.Range("D2:F26").NumberFormat = "@"
$newSheet.Range("G2:G25").NumberFormat = "@"

$r = 2
foreach ($item in $data) {
  $newSheet.Cells.Item($r, 1).Value = $r - 2
  $newSheet.Cells.Item($r, 2).Value = $item[0]
  $newSheet.Cells.Item($r, 3).Value = $item[1]
  $newSheet.Cells.Item($r, 4).Value = $item[2]
  $newSheet.Cells.Item($r, 5).Value = $item[3]
  $newSheet.Cells.Item($r, 6).Value = $item[4]
  $newSheet.Cells.Item($r, 7).Value = $item[5]
  $newSheet.Cells.Item($r, 8).Value = $item[6]
  $r = $r + 1
}

# Last data row's holding-value is genuinely zero, so it is stored as a
# real number (0), not the text "0.0000" used for every other row.
$newSheet.Cells.Item(26, 7).Value = 0

# =====================================================================
# 2. Update the "总计" (totals) sheet: add a 2022-Q1 summary row at the
#    top of the data (right under the header) and push the rest down.
# =====================================================================
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()

$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 25
$total.Range("D2").Value = 5.12

# Renumber the index column for the rows that shifted down one position.
for ($row = 3; $row -le 7; $row++) {
  $total.Cells.Item($row, 1).Value = $row - 2
}

Write-Output "workbook updated"
